$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (A1:J1)
$ws.Range("A1").Value = "No"
$ws.Range("B1").Value = "ID Pelatihan"
$ws.Range("C1").Value = "ID Periode"
$ws.Range("D1").Value = "ID User"
$ws.Range("E1").Value = "Tanggal Mulai"
$ws.Range("F1").Value = "Tanggal Selesai"
$ws.Range("G1").Value = "Lokasi"
$ws.Range("H1").Value = "Quota Peserta"
$ws.Range("I1").Value = "Biaya"
$ws.Range("I1").NumberFormat = "@"
$ws.Range("J1").Value = "Input"
$ws.Range("J1").NumberFormat = "@"

# Clear row 2 data (A2:G2) entirely (content + formatting)
$ws.Range("A2:G2").Clear()
# H2 keeps its (Hyperlink) style but loses its value and the hyperlink itself
$ws.Range("H2").Hyperlinks.Delete()
$ws.Range("H2").ClearContents()

$ws.Columns.Item(2).ColumnWidth = 10.8333333333333

# C1:D1 get a fresh style: General number format, plain Calibri 11 font,
# vertical-center alignment (reset first so the stale Text/Date number format goes away)
$ws.Range("C1:D1").Style = "Normal"
$ws.Range("C1:D1").Font.Name = "Calibri"
$ws.Range("C1:D1").Font.Size = 11
$ws.Range("C1:D1").VerticalAlignment = -4108

# Move the active selection to D6, matching the saved view state
[void]$ws.Range("D6").Select()
